$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift the existing "PAF1/PAF2/PAF3" example block down one row and
#     rename the labels, leaving A1/C1/C2 intro text untouched. ---

# Clear out the old block (rows 4-8) first so nothing stale is left behind
# once things move to their new rows.
$ws.Range("A4:D8").ClearContents()

# Row 4: new section header
$ws.Range("A4").Value = "For exp=8.5"

# Row 5 (was row 4): PAF1 -> "PAF lowci"
$ws.Range("A5").Value = "PAF lowci"
$ws.Range("B5").Formula = "=(1.043879-1)*0.2"

# Row 7 (was row 6): PAF3 -> "PAF3 highci" (written before row 6 below)
$ws.Range("A7").Value = "PAF3 highci"
$ws.Range("B7").Formula = "=(1.167057-1)*0.6"

# Row 6 (was row 5): PAF2 -> "PAF mean"
$ws.Range("A6").Value = "PAF mean"
$ws.Range("B6").Formula = "=(1.103751-1)*0.2"

# Row 9 (was row 8): total
$ws.Range("B9").Formula = "=SUM(B5:B7)"

# --- New block further down the sheet ---
$ws.Range("A12").Value = "Mean rr = 1.118"

$ws.Range("A13").Value = "PAF exp = 8.5"
$ws.Range("B13").Formula = "=(1.04387884019449-1)*0.2"
$ws.Range("D13").Formula = "=(8.85-5)*(1.04387884019449-1)*0.2"

$ws.Range("A14").Value = "PAF exp = 13.85"
$ws.Range("B14").Formula = "=(1.10375071048696-1)*0.2"
$ws.Range("D14").Formula = "=(13.85-5)*(1.10375071048696-1)*0.2"

$ws.Range("A15").Value = "PAF exp = 18.85"
$ws.Range("B15").Formula = "=(1.16705654333744-1)*0.2"
$ws.Range("D15").Formula = "=(18.85-5)*(1.16705654333744-1)*0.6"

# --- Column widths (characters) / selection / view tweaks ---
$ws.Columns.Item(1).ColumnWidth = 11.9
$ws.Columns.Item(2).ColumnWidth = 16

$ws.Range("B12").Select() | Out-Null

# Window got taller by a hair in the source commit too.
$excel.ActiveWindow.Height = 12651 | Out-Null
